# Update the "想去人数" (F column) counts that changed between the two
# data-scrape snapshots (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 56
$ws1.Range("F11").Value = 126
$ws1.Range("F12").Value = 290
$ws1.Range("F15").Value = 684
$ws1.Range("F21").Value = 141
$ws1.Range("F22").Value = 666
$ws1.Range("F27").Value = 871
$ws1.Range("F30").Value = 41
$ws1.Range("F31").Value = 269
$ws1.Range("F33").Value = 15

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 252

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 56
$ws4.Range("F13").Value = 126
$ws4.Range("F14").Value = 290
$ws4.Range("F17").Value = 684
$ws4.Range("F27").Value = 252
$ws4.Range("F28").Value = 252
$ws4.Range("F29").Value = 141
$ws4.Range("F30").Value = 666
$ws4.Range("F35").Value = 871
$ws4.Range("F40").Value = 41
$ws4.Range("F41").Value = 269
$ws4.Range("F45").Value = 15
